$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 140 (pushes existing rows 140.. down to 142..)
$ws.Range("A140:A141").EntireRow.Insert()

# Make sure the date column keeps the same date formatting as the rest of column D
$ws.Range("D140:D141").NumberFormat = $ws.Range("D139").NumberFormat

# New row 140: Choclo / Choclero / Primera, 2023-03-10 (serial 44995)
$ws.Range("A140").Value = 11
$ws.Range("B140").Value = "Vega Monumental Concepción"
$ws.Range("C140").Value = "Bíobío"
$ws.Range("D140").Value = 44995
$ws.Range("E140").Value = 8
$ws.Range("F140").Value = 100112024
$ws.Range("G140").Value = "Choclo"
$ws.Range("H140").Value = "Choclero"
$ws.Range("I140").Value = "Primera"
$ws.Range("J140").Value = 10000
$ws.Range("K140").Value = 350
$ws.Range("L140").Value = 400
$ws.Range("M140").Value = 375
$ws.Range("N140").Value = "$/unidad"
$ws.Range("O140").Value = "Región de O'Higgins"
$ws.Range("P140").Value = 375
$ws.Range("Q140").Value = 1
$ws.Range("R140").Value = "Hortaliza"

# New row 141: Choclo / Choclero / Segunda, 2023-03-10 (serial 44995)
$ws.Range("A141").Value = 11
$ws.Range("B141").Value = "Vega Monumental Concepción"
$ws.Range("C141").Value = "Bíobío"
$ws.Range("D141").Value = 44995
$ws.Range("E141").Value = 8
$ws.Range("F141").Value = 100112024
$ws.Range("G141").Value = "Choclo"
$ws.Range("H141").Value = "Choclero"
$ws.Range("I141").Value = "Segunda"
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 300
$ws.Range("L141").Value = 300
$ws.Range("M141").Value = 300
$ws.Range("N141").Value = "$/unidad"
$ws.Range("O141").Value = "Región de O'Higgins"
$ws.Range("P141").Value = 300
$ws.Range("Q141").Value = 1
$ws.Range("R141").Value = "Hortaliza"
